# Apply attendance_reports sync edit to the "Session Analysis Results" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Row 2 (ANATOMY session 1): reorder "Recorded By" list
$ws.Range("G2").Value = "Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, System, gehanadel@med.asu.edu.eg"

# Row 3 (ANATOMY session 2): reorder "Recorded By" list
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, System"

# Row 4 (ANATOMY session 3): add a recorder and update attendance count
$ws.Range("G4").Value = "gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("H4").Value = "60/251"

# Row 9 (HISTOLOGY session 1): reorder "Recorded By" list
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

# Row 10 / summary row 15: update average attendance percentage.
# Force text format first so Excel keeps these as literal "24.9%" strings
# instead of auto-converting them into a percentage number (which would
# change the stored style and round the displayed value to "24.90%").
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "24.9%"
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "24.9%"

# Row 28 (PHYSIOLOGY session 1): reorder "Recorded By" list
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
